$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3205.6667
$ws.Range("I64").Value = 3120.923
$ws.Range("J64").Value = 3426
$ws.Range("K64").Value = 3120.923
$ws.Range("L64").Value = 3426
$ws.Range("M64").Value = -2872.923
$ws.Range("N64").Value = -3922
$ws.Range("H67").Value = 3205.6667
$ws.Range("I67").Value = 3120.923
$ws.Range("J67").Value = 3426
$ws.Range("K67").Value = 3120.923
$ws.Range("L67").Value = 3426
$ws.Range("M67").Value = -2262.923
$ws.Range("N67").Value = -5142
$ws.Range("H76").Value = 2852874.8
$ws.Range("I76").Value = 3370583.2
$ws.Range("J76").Value = 5477.5
$ws.Range("K76").Value = 3370583.2
$ws.Range("L76").Value = 5477.5
$ws.Range("M76").Value = -3370268.2
$ws.Range("N76").Value = -6107.5
$ws.Range("H79").Value = 2852874.8
$ws.Range("I79").Value = 3370583.2
$ws.Range("J79").Value = 5477.5
$ws.Range("K79").Value = 3370583.2
$ws.Range("L79").Value = 5477.5
$ws.Range("M79").Value = -3369491.2
$ws.Range("N79").Value = -7661.5
$ws.Range("H116").Value = 7835.36
$ws.Range("I116").Value = 11007.692
$ws.Range("J116").Value = 4398.6665
$ws.Range("K116").Value = 11007.692
$ws.Range("L116").Value = 4398.6665
$ws.Range("M116").Value = -7565.691999999999
$ws.Range("N116").Value = -11282.6665
$ws.Range("H137").Value = 1039
$ws.Range("I137").Value = 894.5909
$ws.Range("J137").Value = 1356.7
$ws.Range("K137").Value = 2683.7727
$ws.Range("L137").Value = 4070.1
$ws.Range("M137").Value = -133.7727
$ws.Range("N137").Value = -9170.1
$ws.Range("H139").Value = 69611.42999999999
$ws.Range("J139").Value = 69611.42999999999
$ws.Range("L139").Value = 69611.42999999999
$ws.Range("N139").Value = -79891.42999999999
$ws.Range("H140").Value = 87403
$ws.Range("J140").Value = 87403
$ws.Range("L140").Value = 87403
$ws.Range("N140").Value = -97763

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5654.31
$ws.Range("I32").Value = 3847.0632
$ws.Range("J32").Value = 12453
$ws.Range("K32").Value = 3847.0632
$ws.Range("L32").Value = 12453
$ws.Range("M32").Value = -3560.0632
$ws.Range("N32").Value = -13027
$ws.Range("H61").Value = 2958.0715
$ws.Range("I61").Value = 2944.4595
$ws.Range("K61").Value = 2944.4595
$ws.Range("M61").Value = -2732.4595
$ws.Range("H74").Value = 940.5333000000001
$ws.Range("I74").Value = 550.6667
$ws.Range("K74").Value = 550.6667
$ws.Range("M74").Value = 323.3333
$ws.Range("H77").Value = 940.5333000000001
$ws.Range("I77").Value = 550.6667
$ws.Range("K77").Value = 2753.3335
$ws.Range("M77").Value = 1614.6665
$ws.Range("H136").Value = 2958.0715
$ws.Range("I136").Value = 2944.4595
$ws.Range("K136").Value = 8833.378499999999
$ws.Range("M136").Value = -6283.378499999999
$ws.Range("H139").Value = 64750
$ws.Range("J139").Value = 64750
$ws.Range("L139").Value = 64750
$ws.Range("N139").Value = -75030
$ws.Range("H141").Value = 61852.07
$ws.Range("J141").Value = 61852.07
$ws.Range("L141").Value = 61852.07
$ws.Range("N141").Value = -72212.07000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1657.6538
$ws.Range("I105").Value = 1647.5714
$ws.Range("J105").Value = 1700
$ws.Range("K105").Value = 1647.5714
$ws.Range("L105").Value = 1700
$ws.Range("M105").Value = 99.42859999999996
$ws.Range("N105").Value = -5194
$ws.Range("H138").Value = 59300
$ws.Range("J138").Value = 59300
$ws.Range("L138").Value = 59300
$ws.Range("N138").Value = -69580
$ws.Range("H140").Value = 76113.164
$ws.Range("J140").Value = 76113.164
$ws.Range("L140").Value = 76113.164
$ws.Range("N140").Value = -86473.164

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1148.1515
$ws.Range("I58").Value = 731.8
$ws.Range("J58").Value = 1495.1111
$ws.Range("K58").Value = 731.8
$ws.Range("L58").Value = 1495.1111
$ws.Range("M58").Value = -528.8
$ws.Range("N58").Value = -1901.1111
$ws.Range("H62").Value = 3459.75
$ws.Range("I62").Value = 3587.6667
$ws.Range("J62").Value = 3076
$ws.Range("K62").Value = 3587.6667
$ws.Range("L62").Value = 3076
$ws.Range("M62").Value = -2963.6667
$ws.Range("N62").Value = -4324
$ws.Range("H65").Value = 3459.75
$ws.Range("I65").Value = 3587.6667
$ws.Range("J65").Value = 3076
$ws.Range("K65").Value = 17938.3335
$ws.Range("L65").Value = 15380
$ws.Range("M65").Value = -14818.3335
$ws.Range("N65").Value = -21620
$ws.Range("H136").Value = 1148.1515
$ws.Range("I136").Value = 731.8
$ws.Range("J136").Value = 1495.1111
$ws.Range("K136").Value = 2195.4
$ws.Range("L136").Value = 4485.3333
$ws.Range("M136").Value = 354.6000000000004
$ws.Range("N136").Value = -9585.3333
$ws.Range("H138").Value = 45000
$ws.Range("J138").Value = 45000
$ws.Range("L138").Value = 45000
$ws.Range("N138").Value = -55280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 782.52
$ws.Range("I131").Value = 421.66666
$ws.Range("J131").Value = 831.7273
$ws.Range("K131").Value = 1264.99998
$ws.Range("L131").Value = 2495.1819
$ws.Range("M131").Value = 3775.00002
$ws.Range("N131").Value = -12575.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19404790
$ws.Range("I70").Value = 29612336
$ws.Range("J70").Value = 10450.9
$ws.Range("K70").Value = 29612336
$ws.Range("L70").Value = 10450.9
$ws.Range("M70").Value = -29612066
$ws.Range("N70").Value = -10990.9
$ws.Range("H73").Value = 19404790
$ws.Range("I73").Value = 29612336
$ws.Range("J73").Value = 10450.9
$ws.Range("K73").Value = 29612336
$ws.Range("L73").Value = 10450.9
$ws.Range("M73").Value = -29611400
$ws.Range("N73").Value = -12322.9
$ws.Range("H138").Value = 63306.25
$ws.Range("J138").Value = 63306.25
$ws.Range("L138").Value = 63306.25
$ws.Range("N138").Value = -73586.25
$ws.Range("H140").Value = 99864.5
$ws.Range("J140").Value = 99864.5
$ws.Range("L140").Value = 99864.5
$ws.Range("N140").Value = -110224.5
$ws.Range("H141").Value = 66714.25
$ws.Range("J141").Value = 66714.25
$ws.Range("L141").Value = 66714.25
$ws.Range("N141").Value = -77074.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 52541.547
$ws.Range("J138").Value = 52541.547
$ws.Range("L138").Value = 52541.547
$ws.Range("N138").Value = -62821.547
$ws.Range("H139").Value = 69650
$ws.Range("J139").Value = 69650
$ws.Range("L139").Value = 69650
$ws.Range("N139").Value = -79930

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 54683.332
$ws.Range("J139").Value = 54683.332
$ws.Range("L139").Value = 54683.332
$ws.Range("N139").Value = -64963.332
$ws.Range("H141").Value = 71459.28999999999
$ws.Range("J141").Value = 69702.5
$ws.Range("L141").Value = 69702.5
$ws.Range("N141").Value = -80062.5
